# Applies the "Updated cryptos list" data refresh described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.099.65"
$ws.Range("E2").Value = "  -3.51%  "

$ws.Range("D3").Value = "3.138.58"
$ws.Range("E3").Value = "  -4.91%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'524.34"
$ws.Range("E5").Value = "  -6.02%  "

$ws.Range("D6").Value = "'135.41"
$ws.Range("E6").Value = "  -4.53%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "3.144.48"
$ws.Range("E8").Value = "  -4.77%  "

$ws.Range("E9").Value = "  -5.66%  "

$ws.Range("D10").Value = "'7.25"
$ws.Range("E10").Value = "  -7.46%  "

$ws.Range("E11").Value = "  -6.99%  "

$ws.Range("D12").Value = "'0.386"
$ws.Range("E12").Value = "  -4.79%  "

$ws.Range("D13").Value = "3.672.36"
$ws.Range("E13").Value = "  -5.09%  "

$ws.Range("D14").Value = "'0.127"
$ws.Range("E14").Value = "  -1.98%  "

$ws.Range("D15").Value = "'25.60"
$ws.Range("E15").Value = "  -4.85%  "

$ws.Range("D16").Value = "3.129.44"
$ws.Range("E16").Value = "  -4.92%  "

$ws.Range("D17").Value = "57.992.72"
$ws.Range("E17").Value = "  -3.69%  "

$ws.Range("D18").Value = "'0.0000153"
$ws.Range("E18").Value = "  -7.35%  "

$ws.Range("D19").Value = "'5.84"
$ws.Range("E19").Value = "  -5.74%  "

$ws.Range("D20").Value = "'13.06"
$ws.Range("E20").Value = "  -9.33%  "

$ws.Range("D21").Value = "'8.04"
$ws.Range("E21").Value = "  -6.82%  "

$ws.Range("D22").Value = "'346.08"
$ws.Range("E22").Value = "  -7.58%  "

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").Value = "'68.97"
$ws.Range("E24").Value = "  -7.11%  "

$ws.Range("D25").Value = "'0.506"
$ws.Range("E25").Value = "  -6.49%  "

$ws.Range("D26").Value = "3.254.00"
$ws.Range("E26").Value = "  -5.53%  "

$ws.Range("D27").Value = "'0.169"
$ws.Range("E27").Value = "  -1.76%  "

$ws.Range("D28").Value = "0.0₃0956"
$ws.Range("E28").Value = "  -6.94%  "

$ws.Range("D29").Value = "'0.996"
$ws.Range("E29").Value = "  -0.26%  "

$ws.Range("D30").Value = "'6.84"
$ws.Range("E30").Value = "  -5.41%  "

$ws.Range("D31").Value = "'0.998"
$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("E32").Value = "  -8.10%  "

$ws.Range("E33").Value = "  -9.08%  "

$ws.Range("D34").Value = "'1.25"
$ws.Range("E34").Value = "  -1.21%  "

$ws.Range("D35").Value = "'21.63"
$ws.Range("E35").Value = "  -4.10%  "

$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'4.83"
$ws.Range("E36").Value = "  -5.89%  "

$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "'158.75"
$ws.Range("E37").Value = "  -4.22%  "

$ws.Range("D38").Value = "'6.21"
$ws.Range("E38").Value = "  -7.50%  "

$ws.Range("D39").Value = "'1.39"
$ws.Range("E39").Value = "  -8.99%  "

$ws.Range("D40").Value = "'25.44"
$ws.Range("E40").Value = "  -5.25%  "

$ws.Range("D41").Value = "'0.0694"
$ws.Range("E41").Value = "  -5.93%  "

$ws.Range("D42").Value = "3.170.00"
$ws.Range("E42").Value = "  -4.85%  "

$ws.Range("D43").Value = "'40.27"
$ws.Range("E43").Value = "  -4.02%  "

$ws.Range("D44").Value = "'0.688"
$ws.Range("E44").Value = "  -8.42%  "

$ws.Range("D45").Value = "'1.08"
$ws.Range("E45").Value = "  -3.26%  "

$ws.Range("E46").Value = "  -6.34%  "

$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = "  +0.08%  "

$ws.Range("E48").Value = "  -7.78%  "

$ws.Range("D49").Value = "2.255.68"
$ws.Range("E49").Value = "  -4.76%  "

$ws.Range("E50").Value = "  -4.60%  "

$ws.Range("D51").Value = "'20.34"
$ws.Range("E51").Value = "  -4.31%  "
